$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy formatting (style s="1", used by all existing body rows) from the
# last existing data row down onto the new rows so the new cells reuse the
# existing cellXfs entry instead of minting a new style.
$ws.Rows.Item(712).Copy()
$ws.Range("A713:C739").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new "lab.liquid.*" translation rows (Common: Added liquid section).
$ws.Cells.Item(713,1).Value = "cs"
$ws.Cells.Item(713,2).Value = "lab.liquid.title"
$ws.Cells.Item(713,3).Value = "Liquidy"
$ws.Cells.Item(714,1).Value = "cs"
$ws.Cells.Item(714,2).Value = "lab.liquid.subtitle"
$ws.Cells.Item(714,3).Value = "Tato sekce slouží pro správu liquidů."
$ws.Cells.Item(715,1).Value = "cs"
$ws.Cells.Item(715,2).Value = "lab.liquid.button.create"
$ws.Cells.Item(715,3).Value = "Nový liquid"
$ws.Cells.Item(716,1).Value = "cs"
$ws.Cells.Item(716,2).Value = "lab.liquid.button.list"
$ws.Cells.Item(716,3).Value = "Seznam liquidů"
$ws.Cells.Item(717,1).Value = "cs"
$ws.Cells.Item(717,2).Value = "lab.liquid.list.title"
$ws.Cells.Item(717,3).Value = "Seznam liquidů"
$ws.Cells.Item(718,1).Value = "cs"
$ws.Cells.Item(718,2).Value = "lab.liquid.preview"
$ws.Cells.Item(718,3).Value = "Náhled liquidu"
$ws.Cells.Item(719,1).Value = "cs"
$ws.Cells.Item(719,2).Value = "lab.liquid.button.index"
$ws.Cells.Item(719,3).Value = "Detail liquidu"
$ws.Cells.Item(720,1).Value = "cs"
$ws.Cells.Item(720,2).Value = "lab.liquid.button.edit"
$ws.Cells.Item(720,3).Value = "Editace liquidu"
$ws.Cells.Item(721,1).Value = "cs"
$ws.Cells.Item(721,2).Value = "lab.liquid.button.delete"
$ws.Cells.Item(721,3).Value = "Odstranit liquid"
$ws.Cells.Item(722,1).Value = "cs"
$ws.Cells.Item(722,2).Value = "lab.liquid.button.delete.confirm.title"
$ws.Cells.Item(722,3).Value = "Odstranit liquid"
$ws.Cells.Item(723,1).Value = "cs"
$ws.Cells.Item(723,2).Value = "lab.liquid.button.delete.confirm"
$ws.Cells.Item(723,3).Value = "Tuto akci si prosím velmi rozvažte, poněvadž může mít hluboký dopad na celou aplikaci, hlavně pak statistická data o vapování. Smazání liquidu obecně není doporučeno, proto prosím postupujte s maximální obezřetností, protože není cesty zpět. Vážně."
$ws.Cells.Item(724,1).Value = "cs"
$ws.Cells.Item(724,2).Value = "lab.liquid.button.delete.confirm.ok"
$ws.Cells.Item(724,3).Value = "Odstranit liquid"
$ws.Cells.Item(725,1).Value = "cs"
$ws.Cells.Item(725,2).Value = "lab.liquid.deleted.success"
$ws.Cells.Item(725,3).Value = "Liquid [{{data.name}}] (a možná i hromada jiných dat) byl úspěšně odstraněn."
$ws.Cells.Item(726,1).Value = "cs"
$ws.Cells.Item(726,2).Value = "lab.liquid.table.name"
$ws.Cells.Item(726,3).Value = "Název"
$ws.Cells.Item(727,1).Value = "cs"
$ws.Cells.Item(727,2).Value = "lab.liquid.table.pgvg"
$ws.Cells.Item(727,3).Value = "PG/VG"
$ws.Cells.Item(728,1).Value = "cs"
$ws.Cells.Item(728,2).Value = "lab.liquid.preview.preview.title"
$ws.Cells.Item(728,3).Value = "Detail liquidu"
$ws.Cells.Item(729,1).Value = "cs"
$ws.Cells.Item(729,2).Value = "lab.liquid.preview.preview.subtitle"
$ws.Cells.Item(729,3).Value = "Zde jsou veškeré dostupné informace o vybraném liquidu."
$ws.Cells.Item(730,1).Value = "cs"
$ws.Cells.Item(730,2).Value = "lab.liquid.preview.name"
$ws.Cells.Item(730,3).Value = "Jméno"
$ws.Cells.Item(731,1).Value = "cs"
$ws.Cells.Item(731,2).Value = "lab.liquid.preview.pgvg"
$ws.Cells.Item(731,3).Value = "PG/VG"
$ws.Cells.Item(732,1).Value = "cs"
$ws.Cells.Item(732,2).Value = "lab.liquid.edit.title"
$ws.Cells.Item(732,3).Value = "Editace liquidu"
$ws.Cells.Item(733,1).Value = "cs"
$ws.Cells.Item(733,2).Value = "lab.liquid.edit.subtitle"
$ws.Cells.Item(733,3).Value = "Uprava údajů o vybraném liquidu."
$ws.Cells.Item(734,1).Value = "cs"
$ws.Cells.Item(734,2).Value = "lab.liquid.update.submit"
$ws.Cells.Item(734,3).Value = "Aktualizovat"
$ws.Cells.Item(735,1).Value = "cs"
$ws.Cells.Item(735,2).Value = "lab.liquid.updated.message"
$ws.Cells.Item(735,3).Value = "Liquid [{{data.name}}] byl aktualizován."
$ws.Cells.Item(736,1).Value = "cs"
$ws.Cells.Item(736,2).Value = "lab.liquid.index.title"
$ws.Cells.Item(736,3).Value = "Detail liquidu"
$ws.Cells.Item(737,1).Value = "cs"
$ws.Cells.Item(737,2).Value = "lab.liquid.index.preview.title"
$ws.Cells.Item(737,3).Value = "Detail liquidu"
$ws.Cells.Item(738,1).Value = "cs"
$ws.Cells.Item(738,2).Value = "lab.liquid.index.preview.subtitle"
$ws.Cells.Item(738,3).Value = "Zde jsou veškeré dostupné informace o vybraném liquidu."
$ws.Cells.Item(739,1).Value = "cs"
$ws.Cells.Item(739,2).Value = "lab.liquid.preview.volume"
$ws.Cells.Item(739,3).Value = "Objem"

# Row 723 holds a long wrapped confirmation message; match the taller row height
# used by the other wrapped-text rows in this sheet.
$ws.Rows.Item(723).RowHeight = 39

# Reproduce the view state left behind by the editing session (scrolled down to
# the newly-added rows, with the last new cell selected).
$ws.Activate()
$excel.Goto($ws.Range("A722"), $true)
$ws.Range("B738").Select()
